# Shopping cart page complete module push
# - Clear the Quantity (column D) values from Sheet1 (rows 2-39), leaving the
#   cells empty so the row "spans" shrink back to just column A where D was
#   the last populated column.
# - Update the active selection in the sheet view from C6:C7 to a single
#   cell G4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2:D39").ClearContents()

$ws.Activate()
$ws.Range("G4").Select()
